$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 7 down to the new rows 8 and 9 first
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new data rows 8 and 9, mirroring the style/format of rows 5-7
$ws.Range("A8").Value = 41559
$ws.Range("B8").Value = 0.041666666666666664

$ws.Range("A9").Value = 41560
$ws.Range("B9").Value = 0.1875

# Update selection / view to match the diff
$ws.Range("B4:B9").Select()
